# "For Tim, for augmentation." — add the 14-epoch (56 min, PER) augmented-bin
# run as a third data column (G) alongside the existing 128_bin / 128_bin_times_10
# results on the "new_results" sheet, wire it into the summary table + chart,
# and nudge the active selection as the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("new_results")

# --- new raw-results column (G) on the input table -------------------------
$ws.Range("G3").Value = 91580
$ws.Range("G4").Value = 22895
$ws.Range("G5").Value = 27646
$ws.Range("G6").Value = "14 (56 min, PER)"

# --- accuracy / error figures for the new run -------------------------------
$ws.Range("G8").Value = 0.43049999999999999
$ws.Range("G7").Formula = "=1-G8"
$ws.Range("G7").NumberFormat = "0.00%"

# --- pull it into the summary row used by the chart -------------------------
$ws.Range("C13").Formula = "=G8"

# --- leave the selection where the author left it ---------------------------
$ws.Activate()
$ws.Range("L11").Select()

$wb.Save()
